# Update numeric-looking text values in the "Romania Summary" sheet.
# These cells store plain text (e.g. "20.1") in the shared-string table,
# so a leading apostrophe is used to force Excel to keep them as text
# instead of silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = "'20.06"

$ws.Range("B32").Value = "'16.65"
$ws.Range("D32").Value = "'19.05"

$ws.Range("B34").Value = "'21.88"
$ws.Range("C34").Value = "'44.34"
$ws.Range("D34").Value = "'66.22"

$ws.Range("B36").Value = "'87.05"
$ws.Range("C36").Value = "'12.57"
$ws.Range("D36").Value = "'99.62"

$ws.Range("B40").Value = "'13.33"
$ws.Range("C40").Value = "'41.11"
$ws.Range("D40").Value = "'54.44"
